$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H39").Value = 241.65218
$ws.Range("I39").Value = 136.11111
$ws.Range("J39").Value = 309.5
$ws.Range("K39").Value = 408.33333
$ws.Range("L39").Value = 928.5
$ws.Range("M39").Value = -112.33333
$ws.Range("N39").Value = -1520.5
$ws.Range("H76").Value = 7658.391
$ws.Range("I76").Value = 8595.583000000001
$ws.Range("J76").Value = 6636
$ws.Range("K76").Value = 8595.583000000001
$ws.Range("L76").Value = 6636
$ws.Range("M76").Value = -8280.583000000001
$ws.Range("N76").Value = -7266
$ws.Range("H79").Value = 7658.391
$ws.Range("I79").Value = 8595.583000000001
$ws.Range("J79").Value = 6636
$ws.Range("K79").Value = 8595.583000000001
$ws.Range("L79").Value = 6636
$ws.Range("M79").Value = -7503.583000000001
$ws.Range("N79").Value = -8820
$ws.Range("H116").Value = 6938.7856
$ws.Range("I116").Value = 6714.3
$ws.Range("J116").Value = 7500
$ws.Range("K116").Value = 6714.3
$ws.Range("L116").Value = 7500
$ws.Range("M116").Value = -3272.3
$ws.Range("N116").Value = -14384
$ws.Range("H132").Value = 2078.7407
$ws.Range("I132").Value = 1141.9546
$ws.Range("K132").Value = 3425.8638
$ws.Range("M132").Value = -895.8638000000001
$ws.Range("H137").Value = 174733.72
$ws.Range("I137").Value = 233338.38
$ws.Range("J137").Value = 6733.7334
$ws.Range("K137").Value = 700015.14
$ws.Range("L137").Value = 20201.2002
$ws.Range("M137").Value = -697465.14
$ws.Range("N137").Value = -25301.2002
$ws.Range("H138").Value = 2559.6973
$ws.Range("J138").Value = 2696.2837
$ws.Range("L138").Value = 8088.8511
$ws.Range("N138").Value = -18368.8511

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H33").Value = 14999
$ws.Range("J33").Value = 0
$ws.Range("L33").Value = 0
$ws.Range("N33").ClearContents()
$ws.Range("H122").Value = 2940.1875
$ws.Range("I122").Value = 2022.7222
$ws.Range("J122").Value = 4119.7856
$ws.Range("K122").Value = 6068.1666
$ws.Range("L122").Value = 12359.3568
$ws.Range("M122").Value = -3618.1666
$ws.Range("N122").Value = -17259.3568
$ws.Range("H132").Value = 4227.227
$ws.Range("I132").Value = 2196.3142
$ws.Range("K132").Value = 6588.942599999999
$ws.Range("M132").Value = -4058.942599999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H8").Value = 99.90000000000001
$ws.Range("I8").Value = 99.888885
$ws.Range("J8").Value = 100
$ws.Range("K8").Value = 99.888885
$ws.Range("L8").Value = 100
$ws.Range("M8").Value = 40.111115
$ws.Range("N8").Value = -380
$ws.Range("H94").Value = 2122.0715
$ws.Range("I94").Value = 742
$ws.Range("J94").Value = 2888.7778
$ws.Range("K94").Value = 742
$ws.Range("L94").Value = 2888.7778
$ws.Range("M94").Value = -291
$ws.Range("N94").Value = -3790.7778
$ws.Range("H134").Value = 231816.92
$ws.Range("I134").Value = 1475.5834
$ws.Range("K134").Value = 4426.7502
$ws.Range("M134").Value = -1891.7502

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 6351.5
$ws.Range("I7").Value = 210.3077
$ws.Range("K7").Value = 210.3077
$ws.Range("M7").Value = -97.30770000000001
$ws.Range("H22").Value = 1716.6666
$ws.Range("I22").Value = 1266.6666
$ws.Range("K22").Value = 1266.6666
$ws.Range("M22").Value = -916.6666
$ws.Range("H31").Value = 699320.1
$ws.Range("I31").Value = 14510.083
$ws.Range("K31").Value = 14510.083
$ws.Range("M31").Value = -14215.083
$ws.Range("H34").Value = 699320.1
$ws.Range("I34").Value = 14510.083
$ws.Range("K34").Value = 14510.083
$ws.Range("M34").Value = -14308.083
$ws.Range("H41").Value = 39202
$ws.Range("I41").Value = 3970
$ws.Range("K41").Value = 3970
$ws.Range("M41").Value = -3542
$ws.Range("H50").Value = 0
$ws.Range("I50").Value = 0
$ws.Range("J50").Value = 0
$ws.Range("K50").Value = 0
$ws.Range("L50").Value = 0
$ws.Range("M50").ClearContents()
$ws.Range("N50").ClearContents()
$ws.Range("H51").Value = 75000
$ws.Range("J51").Value = 75000
$ws.Range("L51").Value = 75000
$ws.Range("N51").Value = -76472
$ws.Range("H58").Value = 7967.3335
$ws.Range("I58").Value = 10506.5
$ws.Range("J58").Value = 6697.75
$ws.Range("K58").Value = 10506.5
$ws.Range("L58").Value = 6697.75
$ws.Range("M58").Value = -10303.5
$ws.Range("N58").Value = -7103.75
$ws.Range("H61").Value = 75000
$ws.Range("J61").Value = 75000
$ws.Range("L61").Value = 75000
$ws.Range("N61").Value = -75696
$ws.Range("H99").Value = 2079.111
$ws.Range("I99").Value = 1751.4
$ws.Range("J99").Value = 2488.75
$ws.Range("K99").Value = 1751.4
$ws.Range("L99").Value = 2488.75
$ws.Range("M99").Value = -253.4000000000001
$ws.Range("N99").Value = -5484.75
$ws.Range("H107").Value = 4735
$ws.Range("I107").Value = 3103
$ws.Range("K107").Value = 3103
$ws.Range("M107").Value = -1183
$ws.Range("H126").Value = 2079.111
$ws.Range("I126").Value = 1751.4
$ws.Range("J126").Value = 2488.75
$ws.Range("K126").Value = 5254.200000000001
$ws.Range("L126").Value = 7466.25
$ws.Range("M126").Value = -2784.200000000001
$ws.Range("N126").Value = -12406.25
$ws.Range("H136").Value = 7967.3335
$ws.Range("I136").Value = 10506.5
$ws.Range("J136").Value = 6697.75
$ws.Range("K136").Value = 31519.5
$ws.Range("L136").Value = 20093.25
$ws.Range("M136").Value = -28969.5
$ws.Range("N136").Value = -25193.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H50").Value = 499.29413
$ws.Range("I50").Value = 420
$ws.Range("J50").Value = 509.86667
$ws.Range("K50").Value = 1260
$ws.Range("L50").Value = 1529.60001
$ws.Range("M50").Value = -779
$ws.Range("N50").Value = -2491.60001
$ws.Range("H53").Value = 499.29413
$ws.Range("I53").Value = 420
$ws.Range("J53").Value = 509.86667
$ws.Range("K53").Value = 1260
$ws.Range("L53").Value = 1529.60001
$ws.Range("M53").Value = -779
$ws.Range("N53").Value = -2491.60001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H7").Value = 10000000
$ws.Range("J7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("N7").ClearContents()
$ws.Range("H8").Value = 10000000
$ws.Range("J8").Value = 0
$ws.Range("L8").Value = 0
$ws.Range("N8").ClearContents()
$ws.Range("H33").Value = 20000
$ws.Range("J33").Value = 20000
$ws.Range("L33").Value = 20000
$ws.Range("N33").Value = -20504
$ws.Range("H113").Value = 4542.615
$ws.Range("I113").Value = 5001
$ws.Range("J113").Value = 4504.4165
$ws.Range("K113").Value = 5001
$ws.Range("L113").Value = 4504.4165
$ws.Range("M113").Value = -2831
$ws.Range("N113").Value = -8844.416499999999
$ws.Range("H116").Value = 98626
$ws.Range("J116").Value = 98626
$ws.Range("L116").Value = 98626
$ws.Range("N116").Value = -107804
$ws.Range("H132").Value = 34485370
$ws.Range("I132").Value = 34485370
$ws.Range("K132").Value = 103456110
$ws.Range("M132").Value = -103453580

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 17223.25
$ws.Range("I7").Value = 15430
$ws.Range("K7").Value = 15430
$ws.Range("M7").Value = -15318
$ws.Range("H22").Value = 3096.95
$ws.Range("I22").Value = 2997.0625
$ws.Range("J22").Value = 3496.5
$ws.Range("K22").Value = 2997.0625
$ws.Range("L22").Value = 3496.5
$ws.Range("M22").Value = -2702.0625
$ws.Range("N22").Value = -4086.5
$ws.Range("H27").Value = 3096.95
$ws.Range("I27").Value = 2997.0625
$ws.Range("J27").Value = 3496.5
$ws.Range("K27").Value = 2997.0625
$ws.Range("L27").Value = 3496.5
$ws.Range("M27").Value = -2890.0625
$ws.Range("N27").Value = -3710.5
$ws.Range("H46").Value = 2816.9412
$ws.Range("I46").Value = 2248.875
$ws.Range("K46").Value = 2248.875
$ws.Range("M46").Value = -2060.875
$ws.Range("H55").Value = 76923820
$ws.Range("I55").Value = 125000850
$ws.Range("J55").Value = 578
$ws.Range("K55").Value = 125000850
$ws.Range("L55").Value = 578
$ws.Range("M55").Value = -125000677
$ws.Range("N55").Value = -924
$ws.Range("H61").Value = 1234.9166
$ws.Range("I61").Value = 961.94116
$ws.Range("K61").Value = 961.94116
$ws.Range("M61").Value = -759.94116
$ws.Range("H82").Value = 2982
$ws.Range("I82").Value = 2778.4
$ws.Range("J82").Value = 4000
$ws.Range("K82").Value = 2778.4
$ws.Range("L82").Value = 4000
$ws.Range("M82").Value = -2417.4
$ws.Range("N82").Value = -4722
$ws.Range("H85").Value = 2982
$ws.Range("I85").Value = 2778.4
$ws.Range("J85").Value = 4000
$ws.Range("K85").Value = 2778.4
$ws.Range("L85").Value = 4000
$ws.Range("M85").Value = -1530.4
$ws.Range("N85").Value = -6496
$ws.Range("H113").Value = 1234.9166
$ws.Range("I113").Value = 961.94116
$ws.Range("K113").Value = 961.94116
$ws.Range("M113").Value = 1208.05884
$ws.Range("H118").Value = 138000
$ws.Range("J118").Value = 138000
$ws.Range("L118").Value = 138000
$ws.Range("N118").Value = -141314
$ws.Range("H122").Value = 6261.5
$ws.Range("I122").Value = 5721.0713
$ws.Range("J122").Value = 6681.8335
$ws.Range("K122").Value = 17163.2139
$ws.Range("L122").Value = 20045.5005
$ws.Range("M122").Value = -14713.2139
$ws.Range("N122").Value = -24945.5005
$ws.Range("H126").Value = 17223.25
$ws.Range("I126").Value = 15430
$ws.Range("K126").Value = 46290
$ws.Range("M126").Value = -43820

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H110").Value = 149995
$ws.Range("J110").Value = 149995
$ws.Range("L110").Value = 149995
$ws.Range("N110").Value = -158175
$ws.Range("H116").Value = 131000
$ws.Range("J116").Value = 131000
$ws.Range("L116").Value = 131000
$ws.Range("N116").Value = -140178
$ws.Range("H126").Value = 1632.5
$ws.Range("I126").Value = 1437.1428
$ws.Range("K126").Value = 4311.428400000001
$ws.Range("M126").Value = -1841.428400000001
